$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 780153.75
$ws.Range("J69").Value = 12124.25
$ws.Range("L69").Value = 36372.75
$ws.Range("N69").Value = -38120.75
$ws.Range("H72").Value = 780153.75
$ws.Range("J72").Value = 12124.25
$ws.Range("L72").Value = 109118.25
$ws.Range("N72").Value = -117854.25
$ws.Range("H113").Value = 8498
$ws.Range("I113").Value = 8123.75
$ws.Range("J113").Value = 9995
$ws.Range("K113").Value = 8123.75
$ws.Range("L113").Value = 9995
$ws.Range("M113").Value = -4869.75
$ws.Range("N113").Value = -16503
$ws.Range("H132").Value = 3440.8064
$ws.Range("I132").Value = 2474.64
$ws.Range("K132").Value = 7423.92
$ws.Range("M132").Value = -4893.92
$ws.Range("H137").Value = 2949.2307
$ws.Range("I137").Value = 1743.3914
$ws.Range("J137").Value = 4682.625
$ws.Range("K137").Value = 5230.174199999999
$ws.Range("L137").Value = 14047.875
$ws.Range("M137").Value = -2680.174199999999
$ws.Range("N137").Value = -19147.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4737.206
$ws.Range("I61").Value = 3605.8696
$ws.Range("K61").Value = 3605.8696
$ws.Range("M61").Value = -3393.8696
$ws.Range("H63").Value = 8406.125
$ws.Range("I63").Value = 8312.25
$ws.Range("K63").Value = 8312.25
$ws.Range("M63").Value = -7626.25
$ws.Range("H66").Value = 8406.125
$ws.Range("I66").Value = 8312.25
$ws.Range("K66").Value = 41561.25
$ws.Range("M66").Value = -38129.25
$ws.Range("H74").Value = 4599.7144
$ws.Range("I74").Value = 4231.222
$ws.Range("J74").Value = 5263
$ws.Range("K74").Value = 4231.222
$ws.Range("L74").Value = 5263
$ws.Range("M74").Value = -3357.222
$ws.Range("N74").Value = -7011
$ws.Range("H77").Value = 4599.7144
$ws.Range("I77").Value = 4231.222
$ws.Range("J77").Value = 5263
$ws.Range("K77").Value = 21156.11
$ws.Range("L77").Value = 26315
$ws.Range("M77").Value = -16788.11
$ws.Range("N77").Value = -35051
$ws.Range("H132").Value = 9132.085999999999
$ws.Range("I132").Value = 10605.241
$ws.Range("K132").Value = 31815.723
$ws.Range("M132").Value = -29285.723
$ws.Range("H136").Value = 4737.206
$ws.Range("I136").Value = 3605.8696
$ws.Range("K136").Value = 10817.6088
$ws.Range("M136").Value = -8267.6088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4227.5386
$ws.Range("I86").Value = 4751.6665
$ws.Range("J86").Value = 3778.2856
$ws.Range("K86").Value = 4751.6665
$ws.Range("L86").Value = 3778.2856
$ws.Range("M86").Value = -3628.6665
$ws.Range("N86").Value = -6024.2856
$ws.Range("H89").Value = 4227.5386
$ws.Range("I89").Value = 4751.6665
$ws.Range("J89").Value = 3778.2856
$ws.Range("K89").Value = 23758.3325
$ws.Range("L89").Value = 18891.428
$ws.Range("M89").Value = -18142.3325
$ws.Range("N89").Value = -30123.428
$ws.Range("H134").Value = 7347.2856
$ws.Range("I134").Value = 4714
$ws.Range("K134").Value = 14142
$ws.Range("M134").Value = -11607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9254.77
$ws.Range("J62").Value = 8802.5
$ws.Range("L62").Value = 8802.5
$ws.Range("N62").Value = -10050.5
$ws.Range("H65").Value = 9254.77
$ws.Range("J65").Value = 8802.5
$ws.Range("L65").Value = 44012.5
$ws.Range("N65").Value = -50252.5
$ws.Range("H99").Value = 4258.49
$ws.Range("I99").Value = 3560.8684
$ws.Range("J99").Value = 6297.6924
$ws.Range("K99").Value = 3560.8684
$ws.Range("L99").Value = 6297.6924
$ws.Range("M99").Value = -2062.8684
$ws.Range("N99").Value = -9293.6924
$ws.Range("H126").Value = 4258.49
$ws.Range("I126").Value = 3560.8684
$ws.Range("J126").Value = 6297.6924
$ws.Range("K126").Value = 10682.6052
$ws.Range("L126").Value = 18893.0772
$ws.Range("M126").Value = -8212.6052
$ws.Range("N126").Value = -23833.0772
$ws.Range("H134").Value = 3330.4285
$ws.Range("I134").Value = 1268.1428
$ws.Range("J134").Value = 9517.286
$ws.Range("K134").Value = 3804.4284
$ws.Range("L134").Value = 28551.858
$ws.Range("M134").Value = -1269.4284
$ws.Range("N134").Value = -33621.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52.526318
$ws.Range("I2").Value = 54.9375
$ws.Range("J2").Value = 39.666668
$ws.Range("K2").Value = 329.625
$ws.Range("L2").Value = 238.000008
$ws.Range("M2").Value = -216.625
$ws.Range("N2").Value = -464.000008
$ws.Range("H42").Value = 5466
$ws.Range("J42").Value = 5899
$ws.Range("L42").Value = 17697
$ws.Range("N42").Value = -18765
$ws.Range("H74").Value = 9343.333000000001
$ws.Range("J74").Value = 12015
$ws.Range("L74").Value = 36045
$ws.Range("N74").Value = -38167
$ws.Range("H77").Value = 9343.333000000001
$ws.Range("J77").Value = 12015
$ws.Range("L77").Value = 108135
$ws.Range("N77").Value = -118743
$ws.Range("H82").Value = 8308
$ws.Range("I82").Value = 7462
$ws.Range("K82").Value = 22386
$ws.Range("M82").Value = -21980
$ws.Range("H85").Value = 8308
$ws.Range("I85").Value = 7462
$ws.Range("K85").Value = 22386
$ws.Range("M85").Value = -20982
$ws.Range("H92").Value = 1036.8334
$ws.Range("I92").Value = 160
$ws.Range("K92").Value = 480
$ws.Range("M92").Value = 768
$ws.Range("H136").Value = 2629.7896
$ws.Range("I136").Value = 2248.875
$ws.Range("K136").Value = 6746.625
$ws.Range("M136").Value = -1646.625
$ws.Range("H140").Value = 1610.8125
$ws.Range("I140").Value = 1126.6428
$ws.Range("K140").Value = 3379.9284
$ws.Range("M140").Value = 1800.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8273.272000000001
$ws.Range("I132").Value = 8502.25
$ws.Range("J132").Value = 8142.4287
$ws.Range("K132").Value = 25506.75
$ws.Range("L132").Value = 24427.2861
$ws.Range("M132").Value = -22976.75
$ws.Range("N132").Value = -29487.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66750584
$ws.Range("I7").Value = 125152344
$ws.Range("K7").Value = 125152344
$ws.Range("M7").Value = -125152232
$ws.Range("H40").Value = 102023.32
$ws.Range("I40").Value = 123478.055
$ws.Range("K40").Value = 123478.055
$ws.Range("M40").Value = -123342.055
$ws.Range("H82").Value = 41668670
$ws.Range("I82").Value = 2287.1875
$ws.Range("J82").Value = 125001440
$ws.Range("K82").Value = 2287.1875
$ws.Range("L82").Value = 125001440
$ws.Range("M82").Value = -1926.1875
$ws.Range("N82").Value = -125002162
$ws.Range("H85").Value = 41668670
$ws.Range("I85").Value = 2287.1875
$ws.Range("J85").Value = 125001440
$ws.Range("K85").Value = 2287.1875
$ws.Range("L85").Value = 125001440
$ws.Range("M85").Value = -1039.1875
$ws.Range("N85").Value = -125003936
$ws.Range("H126").Value = 66750584
$ws.Range("I126").Value = 125152344
$ws.Range("K126").Value = 375457032
$ws.Range("M126").Value = -375454562

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 92851.52
$ws.Range("I62").Value = 503648
$ws.Range("J62").Value = 6368.0527
$ws.Range("K62").Value = 503648
$ws.Range("L62").Value = 6368.0527
$ws.Range("M62").Value = -503024
$ws.Range("N62").Value = -7616.0527
$ws.Range("H65").Value = 92851.52
$ws.Range("I65").Value = 503648
$ws.Range("J65").Value = 6368.0527
$ws.Range("K65").Value = 2518240
$ws.Range("L65").Value = 31840.2635
$ws.Range("M65").Value = -2515120
$ws.Range("N65").Value = -38080.2635
$ws.Range("H96").Value = 2614.2
$ws.Range("I96").Value = 2719.182
$ws.Range("J96").Value = 2325.5
$ws.Range("K96").Value = 2719.182
$ws.Range("L96").Value = 2325.5
$ws.Range("M96").Value = -1346.182
$ws.Range("N96").Value = -5071.5
$ws.Range("H100").Value = 2960237
$ws.Range("I100").Value = 3497704.5
$ws.Range("J100").Value = 4166.5
$ws.Range("K100").Value = 6995409
$ws.Range("L100").Value = 8333
$ws.Range("M100").Value = -6994868
$ws.Range("N100").Value = -9415
$ws.Range("H126").Value = 25701.21
$ws.Range("I126").Value = 33640.23
$ws.Range("K126").Value = 100920.69
$ws.Range("M126").Value = -98450.69
